# Apply cryptos list update (prices / volume changes + EnergySwap/WEMIXToken row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "49.751.26"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.94%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.664.18"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +7.26%  "
$ws.Range("E4").Value = "  +0.18%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "113.02"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +8.29%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "325.47"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.50%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.526"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.48%  "
$ws.Range("E8").Value = "  +0.12%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.550"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.89%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.48"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.81%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.03"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0819"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("E13").Value = "  +0.09%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.32"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.084.96"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.39%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.660.98"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +7.16%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.868"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.38%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "49.781.25"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.27%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.06"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.27%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.74"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.53%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.91"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.0₃0956"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "71.48"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.27%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "274.21"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("E25").Value = "  +3.41%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.70"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("E27").Value = "  +0.08%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.14"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.06%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.87"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +4.74%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "50.12"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.19%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.46"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.54%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "19.37"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.73%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0803"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("E36").Value = "  -0.06%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +11.58%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.56%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.12"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +9.38%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "125.29"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +5.31%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.21"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "22.11"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +4.23%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0316"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.41%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.110.28"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.35%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("E47").Value = "  +9.43%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.21"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +5.76%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.26%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.91%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "59.36"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.16%  "
